$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 460
$ws.Range("I5").Value = 412
$ws.Range("K5").Value = 412
$ws.Range("M5").Value = -297
$ws.Range("H17").Value = 985.4643
$ws.Range("J17").Value = 985.4643
$ws.Range("L17").Value = 2956.3929
$ws.Range("N17").Value = -3292.3929
$ws.Range("H70").Value = 12094.889
$ws.Range("I70").Value = 650
$ws.Range("K70").Value = 1950
$ws.Range("M70").Value = -1680
$ws.Range("H73").Value = 12094.889
$ws.Range("I73").Value = 650
$ws.Range("K73").Value = 1950
$ws.Range("M73").Value = -1014
$ws.Range("H80").Value = 1007.36365
$ws.Range("I80").Value = 1359.6923
$ws.Range("J80").Value = 498.44446
$ws.Range("K80").Value = 4079.0769
$ws.Range("L80").Value = 1495.33338
$ws.Range("M80").Value = -3081.0769
$ws.Range("N80").Value = -3491.33338
$ws.Range("H83").Value = 1007.36365
$ws.Range("I83").Value = 1359.6923
$ws.Range("J83").Value = 498.44446
$ws.Range("K83").Value = 12237.2307
$ws.Range("L83").Value = 4486.00014
$ws.Range("M83").Value = -7245.2307
$ws.Range("N83").Value = -14470.00014
$ws.Range("H132").Value = 1151.4
$ws.Range("I132").Value = 984.7778
$ws.Range("J132").Value = 2651
$ws.Range("K132").Value = 2954.3334
$ws.Range("L132").Value = 7953
$ws.Range("M132").Value = -424.3334
$ws.Range("N132").Value = -13013
$ws.Range("H135").Value = 534.9048
$ws.Range("I135").Value = 518.1667
$ws.Range("J135").Value = 635.3333
$ws.Range("K135").Value = 4663.5003
$ws.Range("L135").Value = 5717.9997
$ws.Range("M135").Value = -2128.5003
$ws.Range("N135").Value = -10787.9997
$ws.Range("H138").Value = 2914.5715
$ws.Range("I138").Value = 2586.7407
$ws.Range("J138").Value = 3316.9092
$ws.Range("K138").Value = 7760.222099999999
$ws.Range("L138").Value = 9950.7276
$ws.Range("M138").Value = -2620.222099999999
$ws.Range("N138").Value = -20230.7276

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 342.375
$ws.Range("I5").Value = 194.83333
$ws.Range("K5").Value = 194.83333
$ws.Range("M5").Value = -82.83332999999999
$ws.Range("H32").Value = 3843.8547
$ws.Range("I32").Value = 3159.7754
$ws.Range("K32").Value = 3159.7754
$ws.Range("M32").Value = -2872.7754
$ws.Range("H53").Value = 15000
$ws.Range("I53").Value = 10000
$ws.Range("K53").Value = 10000
$ws.Range("M53").Value = -9318
$ws.Range("H61").Value = 6136.8335
$ws.Range("I61").Value = 3594.75
$ws.Range("K61").Value = 3594.75
$ws.Range("M61").Value = -3382.75
$ws.Range("H74").Value = 1351.258
$ws.Range("I74").Value = 1085.3914
$ws.Range("J74").Value = 2115.625
$ws.Range("K74").Value = 1085.3914
$ws.Range("L74").Value = 2115.625
$ws.Range("M74").Value = -211.3914
$ws.Range("N74").Value = -3863.625
$ws.Range("H77").Value = 1351.258
$ws.Range("I77").Value = 1085.3914
$ws.Range("J77").Value = 2115.625
$ws.Range("K77").Value = 5426.957
$ws.Range("L77").Value = 10578.125
$ws.Range("M77").Value = -1058.957
$ws.Range("N77").Value = -19314.125
$ws.Range("H136").Value = 6136.8335
$ws.Range("I136").Value = 3594.75
$ws.Range("K136").Value = 10784.25
$ws.Range("M136").Value = -8234.25

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 342.375
$ws.Range("I4").Value = 194.83333
$ws.Range("K4").Value = 194.83333
$ws.Range("M4").Value = -79.83332999999999
$ws.Range("H107").Value = 1527.591
$ws.Range("J107").Value = 1879.8
$ws.Range("L107").Value = 1879.8
$ws.Range("N107").Value = -5719.8
$ws.Range("H57").Value = 49800
$ws.Range("J57").Value = 49800
$ws.Range("L57").Value = 49800
$ws.Range("N57").Value = -51240
$ws.Range("H136").Value = 49800
$ws.Range("J136").Value = 49800
$ws.Range("L136").Value = 49800
$ws.Range("N136").Value = -60000

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 892
$ws.Range("I16").Value = 849.4
$ws.Range("K16").Value = 849.4
$ws.Range("M16").Value = -562.4
$ws.Range("H31").Value = 1846.0416
$ws.Range("I31").Value = 1229.0834
$ws.Range("K31").Value = 1229.0834
$ws.Range("M31").Value = -934.0834
$ws.Range("H34").Value = 1846.0416
$ws.Range("I34").Value = 1229.0834
$ws.Range("K34").Value = 1229.0834
$ws.Range("M34").Value = -1027.0834
$ws.Range("H58").Value = 2290793.2
$ws.Range("J58").Value = 3123.889
$ws.Range("L58").Value = 3123.889
$ws.Range("N58").Value = -3529.889
$ws.Range("H113").Value = 892
$ws.Range("I113").Value = 849.4
$ws.Range("K113").Value = 849.4
$ws.Range("M113").Value = 1320.6
$ws.Range("H134").Value = 1350.0454
$ws.Range("I134").Value = 1385.05
$ws.Range("J134").Value = 1000
$ws.Range("K134").Value = 4155.15
$ws.Range("L134").Value = 3000
$ws.Range("M134").Value = -1620.15
$ws.Range("N134").Value = -8070
$ws.Range("H136").Value = 2290793.2
$ws.Range("J136").Value = 3123.889
$ws.Range("M136").Value = -13046536.5
$ws.Range("N136").Value = -14471.667

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 502000
$ws.Range("I22").Value = 501500
$ws.Range("J22").Value = 502500
$ws.Range("K22").Value = 1504500
$ws.Range("L22").Value = 1507500
$ws.Range("M22").Value = -1504331
$ws.Range("N22").Value = -1507838
$ws.Range("H27").Value = 502000
$ws.Range("I27").Value = 501500
$ws.Range("J27").Value = 502500
$ws.Range("K27").Value = 1504500
$ws.Range("L27").Value = 1507500
$ws.Range("M27").Value = -1504398
$ws.Range("N27").Value = -1507704
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2830479.5
$ws.Range("I126").Value = 2926820.5
$ws.Range("K126").Value = 8780461.5
$ws.Range("M126").Value = -8777991.5

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1610.1052
$ws.Range("I22").Value = 1436.2
$ws.Range("J22").Value = 1803.3334
$ws.Range("K22").Value = 1436.2
$ws.Range("L22").Value = 1803.3334
$ws.Range("M22").Value = -1141.2
$ws.Range("N22").Value = -2393.3334
$ws.Range("H27").Value = 1610.1052
$ws.Range("I27").Value = 1436.2
$ws.Range("J27").Value = 1803.3334
$ws.Range("K27").Value = 1436.2
$ws.Range("L27").Value = 1803.3334
$ws.Range("M27").Value = -1329.2
$ws.Range("N27").Value = -2017.3334
$ws.Range("H98").Value = 99999
$ws.Range("J98").Value = 99999
$ws.Range("L98").Value = 99999
$ws.Range("N98").Value = -105989

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1333.4546
$ws.Range("I132").Value = 1070.3684
$ws.Range("J132").Value = 2999.6667
$ws.Range("K132").Value = 3211.1052
$ws.Range("L132").Value = 8999.000100000001
$ws.Range("M132").Value = -681.1052
$ws.Range("N132").Value = -14059.0001

Write-Output "Applied all Tonberry_Profits updates"
